# Updated cryptos list values (Price column D, Volume(1h) column E)
# Each entry: row number, new D value (or $null to leave unchanged), new E value (or $null to leave unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "27.149.42";  E = "  +0.41%  " },
    @{ Row = 3;  D = "1.679.27";   E = "  -0.04%  " },
    @{ Row = 4;  D = $null;        E = "  +0.11%  " },
    @{ Row = 5;  D = "214.49";     E = "  -0.67%  " },
    @{ Row = 6;  D = "0.517";      E = $null },
    @{ Row = 7;  D = $null;        E = "  +0.09%  " },
    @{ Row = 8;  D = "22.78";      E = "  +6.29%  " },
    @{ Row = 9;  D = "0.261";      E = "  +2.57%  " },
    @{ Row = 10; D = $null;        E = "  -0.44%  " },
    @{ Row = 11; D = "0.0890";     E = "  +0.19%  " },
    @{ Row = 12; D = $null;        E = "  -0.05%  " },
    @{ Row = 13; D = "1.684.35";   E = "  +0.08%  " },
    @{ Row = 14; D = $null;        E = "  +2.19%  " },
    @{ Row = 15; D = $null;        E = "  +4.00%  " },
    @{ Row = 16; D = "66.58";      E = "  +0.26%  " },
    @{ Row = 17; D = "27.122.66";  E = "  +0.31%  " },
    @{ Row = 18; D = "234.59";     E = "  -0.49%  " },
    @{ Row = 19; D = "7.87";       E = "  -3.94%  " },
    @{ Row = 20; D = $null;        E = "  +0.42%  " },
    @{ Row = 21; D = $null;        E = "  +0.08%  " },
    @{ Row = 22; D = $null;        E = "  +1.55%  " },
    @{ Row = 23; D = $null;        E = "  +2.80%  " },
    @{ Row = 25; D = "148.92";     E = "  +1.49%  " },
    @{ Row = 26; D = "7.44";       E = "  +2.36%  " },
    @{ Row = 27; D = "16.35";      E = "  -0.70%  " },
    @{ Row = 29; D = $null;        E = "  -0.04%  " },
    @{ Row = 30; D = $null;        E = "  +0.52%  " },
    @{ Row = 31; D = $null;        E = "  -0.41%  " },
    @{ Row = 32; D = $null;        E = "  -0.14%  " },
    @{ Row = 33; D = "1.540.61";   E = "  +0.16%  " },
    @{ Row = 34; D = $null;        E = "  +0.81%  " },
    @{ Row = 35; D = $null;        E = "  -4.20%  " },
    @{ Row = 36; D = "0.608";      E = "  +3.11%  " },
    @{ Row = 37; D = "0.941";      E = "  +2.78%  " },
    @{ Row = 38; D = $null;        E = "  -0.06%  " },
    @{ Row = 39; D = $null;        E = "  -0.93%  " },
    @{ Row = 40; D = $null;        E = "  +2.24%  " },
    @{ Row = 41; D = $null;        E = "  +2.44%  " },
    @{ Row = 42; D = $null;        E = "  +3.45%  " },
    @{ Row = 43; D = $null;        E = "  +0.12%  " },
    @{ Row = 44; D = $null;        E = "  -0.51%  " },
    @{ Row = 45; D = "1.824.97";   E = "  +0.12%  " },
    @{ Row = 46; D = "0.779";      E = "  +0.06%  " },
    @{ Row = 47; D = "89.72";      E = "  -0.84%  " },
    @{ Row = 48; D = $null;        E = "  +6.51%  " },
    @{ Row = 49; D = $null;        E = "  +2.60%  " },
    @{ Row = 50; D = "8.22";       E = "  +2.38%  " },
    @{ Row = 51; D = $null;        E = "  -0.43%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force the Price column to remain text (it holds dotted
        # pseudo-numbers like "27.149.42" as well as plain decimals like
        # "214.49" that Excel would otherwise silently convert to a real
        # number / mangle via floating point). Save + restore the cell's
        # original style so no formatting/number-format change leaks in.
        $cell = $ws.Cells.Item($u.Row, 4)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = $origStyle
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
